$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.195.02"
$ws.Range("E2").Value = "  +4.66%  "
$ws.Range("D3").Value = "2.225.91"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "83.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.68%  "
$ws.Range("E7").Value = "  +3.19%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "2.561.67"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("D16").Value = "2.209.19"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "44.032.65"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000104"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("E22").Value = "  +10.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0895"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("E34").Value = "  +5.12%  "
$ws.Range("E35").Value = "  +10.12%  "
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0367"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.35%  "
$ws.Range("E38").Value = "  +7.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.56%  "
$ws.Range("E41").Value = "  +4.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.48%  "
$ws.Range("E44").Value = "  +3.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("E51").Value = "  +4.13%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0989"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +30.57%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.43%  "
